$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number-format on price cells whose new values would
# otherwise be auto-coerced to numbers by Excel (single-dot decimals).
$textCells = @("D5","D8","D9","D10","D11","D14","D15","D16","D19","D21","D22","D23","D24","D25","D27","D28","D29","D30","D32","D33","D36","D38","D39","D41","D44","D45","D46","D47","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values
$ws.Range("D2").Value = "26.944.15"
$ws.Range("E2").Value = "  +0.35%  "
$ws.Range("D3").Value = "1.640.23"
$ws.Range("E3").Value = "  -0.06%  "
$ws.Range("E4").Value = "  -0.54%  "
$ws.Range("D5").Value = "217.76"
$ws.Range("E5").Value = "  -0.03%  "
$ws.Range("E6").Value = "  +1.81%  "
$ws.Range("E7").Value = "  -0.50%  "
$ws.Range("D8").Value = "0.254"
$ws.Range("E8").Value = "  +1.50%  "
$ws.Range("D9").Value = "0.0624"
$ws.Range("E9").Value = "  +0.29%  "
$ws.Range("D10").Value = "19.94"
$ws.Range("E10").Value = "  +3.71%  "
$ws.Range("D11").Value = "0.0845"
$ws.Range("E11").Value = "  +0.04%  "
$ws.Range("D12").Value = "1.868.34"
$ws.Range("E12").Value = "  -0.15%  "
$ws.Range("D13").Value = "1.630.70"
$ws.Range("E13").Value = "  -1.05%  "
$ws.Range("D14").Value = "4.11"
$ws.Range("E14").Value = "  -0.98%  "
$ws.Range("D15").Value = "0.533"
$ws.Range("E15").Value = "  +1.30%  "
$ws.Range("D16").Value = "67.24"
$ws.Range("E16").Value = "  +3.27%  "
$ws.Range("D17").Value = "26.921.53"
$ws.Range("E17").Value = "  +0.24%  "
$ws.Range("E18").Value = "  +0.10%  "
$ws.Range("D19").Value = "219.42"
$ws.Range("E19").Value = "  +2.04%  "
$ws.Range("E20").Value = "  -0.47%  "
$ws.Range("D21").Value = "6.75"
$ws.Range("E21").Value = "  +2.30%  "
$ws.Range("D22").Value = "4.41"
$ws.Range("E22").Value = "  +1.16%  "
$ws.Range("D23").Value = "2.43"
$ws.Range("E23").Value = "  +2.10%  "
$ws.Range("D24").Value = "9.20"
$ws.Range("E24").Value = "  +0.27%  "
$ws.Range("D25").Value = "147.34"
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("E26").Value = "  -0.52%  "
$ws.Range("D27").Value = "0.119"
$ws.Range("E27").Value = "  +0.74%  "
$ws.Range("D28").Value = "7.27"
$ws.Range("E28").Value = "  +1.16%  "
$ws.Range("D29").Value = "15.75"
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("D30").Value = "0.0504"
$ws.Range("E30").Value = "  -0.65%  "
$ws.Range("E31").Value = "  -0.84%  "
$ws.Range("D32").Value = "3.34"
$ws.Range("E32").Value = "  -0.53%  "
$ws.Range("D33").Value = "3.01"
$ws.Range("E33").Value = "  +0.57%  "
$ws.Range("E34").Value = "  +1.08%  "
$ws.Range("D35").Value = "1.266.55"
$ws.Range("E35").Value = "  -0.91%  "
$ws.Range("D36").Value = "2.44"
$ws.Range("E36").Value = "  +0.00%  "
$ws.Range("E37").Value = "  +2.77%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").Value = "0.537"
$ws.Range("E38").Value = "  +1.42%  "
$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").Value = "0.840"
$ws.Range("E39").Value = "  +2.53%  "
$ws.Range("E40").Value = "  -0.47%  "
$ws.Range("D41").Value = "0.807"
$ws.Range("E41").Value = "  +0.50%  "
$ws.Range("E42").Value = "  +0.99%  "
$ws.Range("D43").Value = "1.778.99"
$ws.Range("E43").Value = "  -0.21%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "62.42"
$ws.Range("E44").Value = "  +2.32%  "
$ws.Range("B45").Value = "MXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D45").Value = "2.12"
$ws.Range("E45").Value = "  +1.07%  "
$ws.Range("D46").Value = "92.31"
$ws.Range("E46").Value = "  -0.20%  "
$ws.Range("D47").Value = "1.62"
$ws.Range("E47").Value = "  +1.27%  "
$ws.Range("E48").Value = "  +15.08%  "
$ws.Range("E49").Value = "  -0.69%  "
$ws.Range("D50").Value = "7.68"
$ws.Range("E50").Value = "  +1.24%  "
$ws.Range("D51").Value = "0.0963"
$ws.Range("E51").Value = "  -0.42%  "

# Restore default (Normal) style on the cells we forced to text so
# the style table matches the un-styled original cells.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
